$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A14").Value = "EasyBank"
$ws.Range("B14").Value = "https://www.linkedin.com/posts/easybank-official_summerinternships-finance-marketing-activity-7200065600138682368-m-zZ?utm_source=share&utm_medium=member_desktop"
$ws.Range("C14").Value = "IT"
[void]$ws.Range("C14").Select()
